$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) column F values need to be reset to 0
# for rows 2-5 on both the "展览" (Exhibition) sheet and the "全部类型"
# (All Types) sheet, which mirror the same data.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2:F5").Value = 0
}
